# Trade #77 closed at 2026-02-17 21:13:32 - unknown UNKNOWN +0.000%
#
# This script replays the bookkeeping update that happens when the live
# trading bot closes one MarketMaking trade (the prior Trade #105 / row 106
# in "All Trades" & row 73 in "MarketMaking") and opens a brand-new one
# (Trade #138, appended as a new row in both sheets), then rolls the
# aggregate stats on the "Summary" and "Strategy Status" sheets forward.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet - top level P&L / trade-count roll-up
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.14          # Current Capital
$summary.Range("B4").Value = 0.93             # Total P&L $
$summary.Range("B6").Value = 105              # Total Trades
$summary.Range("B8").Value = 41               # Losing Trades
$summary.Range("B9").Value = 46.67            # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.14            # Capital
$status.Range("D5").Value = 72                # Trades
$status.Range("E5").Value = 0.82              # P&L $
$status.Range("F5").Value = 1.14              # P&L %
$status.Range("G5").Value = 48.61             # Win Rate %

# ---------------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# 3a) Close the previously-open Trade #105 (sheet row 106)
$allTrades.Cells.Item(106, 7).Value = 0.02            # Exit Price
$allTrades.Cells.Item(106, 8).Value = "CLOSED"        # Status
$allTrades.Cells.Item(106, 9).Value = -33.3333        # P&L %
$allTrades.Cells.Item(106, 10).Value = -0.01          # P&L $
$allTrades.Cells.Item(106, 11).Value = 101.14         # Capital After
$allTrades.Cells.Item(106, 12).Value = "early_exit"   # Exit Reason
$allTrades.Cells.Item(106, 13).Value = 0.13           # Duration (min)

# 3b) Append the newly-opened Trade #138 as row 139
$r = 139
$allTrades.Cells.Item($r, 1).Value = 138
$allTrades.Cells.Item($r, 2).NumberFormat = "@"
$allTrades.Cells.Item($r, 2).Value = "2026-02-17"
$allTrades.Cells.Item($r, 3).NumberFormat = "@"
$allTrades.Cells.Item($r, 3).Value = "21:13:25"
$allTrades.Cells.Item($r, 4).Value = "MarketMaking"
$allTrades.Cells.Item($r, 5).Value = "DOWN"
$allTrades.Cells.Item($r, 6).Value = 0.03
$allTrades.Cells.Item($r, 7).NumberFormat = "@"
$allTrades.Cells.Item($r, 7).Value = ""
$allTrades.Cells.Item($r, 8).Value = "OPEN"
$allTrades.Cells.Item($r, 9).Value = 0
$allTrades.Cells.Item($r, 10).Value = 0
$allTrades.Cells.Item($r, 11).Value = 101.1496151053151
$allTrades.Cells.Item($r, 12).NumberFormat = "@"
$allTrades.Cells.Item($r, 12).Value = ""
$allTrades.Cells.Item($r, 13).Value = 0
$allTrades.Cells.Item($r, 14).Value = 0
$allTrades.Cells.Item($r, 15).Value = 0
$allTrades.Cells.Item($r, 16).Value = 0.6
$allTrades.Cells.Item($r, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# 4) MarketMaking sheet (strategy-specific trade log)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# 4a) Close the previously-open Trade #105 (sheet row 73)
$mm.Cells.Item(73, 7).Value = 0.02                 # Exit Price
$mm.Cells.Item(73, 8).Value = "CLOSED"             # Status
$mm.Cells.Item(73, 9).Value = -33.3333             # P&L %
$mm.Cells.Item(73, 10).Value = -0.01               # P&L $
$mm.Cells.Item(73, 11).Value = 101.14              # Capital After
$mm.Cells.Item(73, 16).NumberFormat = "@"
$mm.Cells.Item(73, 16).Value = "early_exit"        # Exit Reason
$mm.Cells.Item(73, 17).Value = 0.13                # Duration (min)

# 4b) Append the newly-opened Trade #138 as row 106
$r2 = 106
$mm.Cells.Item($r2, 1).Value = 138
$mm.Cells.Item($r2, 2).NumberFormat = "@"
$mm.Cells.Item($r2, 2).Value = "2026-02-17"
$mm.Cells.Item($r2, 3).NumberFormat = "@"
$mm.Cells.Item($r2, 3).Value = "21:13:25"
$mm.Cells.Item($r2, 4).Value = "MarketMaking"
$mm.Cells.Item($r2, 5).Value = "DOWN"
$mm.Cells.Item($r2, 6).Value = 0.03
$mm.Cells.Item($r2, 7).NumberFormat = "@"
$mm.Cells.Item($r2, 7).Value = ""
$mm.Cells.Item($r2, 8).Value = "OPEN"
$mm.Cells.Item($r2, 9).Value = 0
$mm.Cells.Item($r2, 10).Value = 0
$mm.Cells.Item($r2, 11).Value = 101.1496151053151
$mm.Cells.Item($r2, 12).Value = 0
$mm.Cells.Item($r2, 13).Value = 0
$mm.Cells.Item($r2, 14).Value = 0.6
$mm.Cells.Item($r2, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item($r2, 16).NumberFormat = "@"
$mm.Cells.Item($r2, 16).Value = ""
$mm.Cells.Item($r2, 17).Value = 0
